# Auto-generated edit script: updates market-price-derived columns (H-N)
# on the per-job leve-profit sheets, per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 34950
$ws.Range("J3").Value = 34950
$ws.Range("L3").Value = 34950
$ws.Range("N3").Value = -35178
$ws.Range("H15").Value = 11365048
$ws.Range("I15").Value = 11365048
$ws.Range("K15").Value = 34095144
$ws.Range("M15").Value = -34094975
$ws.Range("H17").Value = 1622.2222
$ws.Range("J17").Value = 1622.2222
$ws.Range("L17").Value = 4866.6666
$ws.Range("N17").Value = -5202.6666
$ws.Range("H33").Value = 726.61536
$ws.Range("I33").Value = 726.61536
$ws.Range("K33").Value = 726.61536
$ws.Range("M33").Value = -497.61536
$ws.Range("H48").Value = 10000
$ws.Range("J48").Value = 10000
$ws.Range("L48").Value = 30000
$ws.Range("N48").Value = -30584
$ws.Range("H56").Value = 10000
$ws.Range("J56").Value = 10000
$ws.Range("L56").Value = 30000
$ws.Range("N56").Value = -31068
$ws.Range("H74").Value = 4759.533
$ws.Range("I74").Value = 4356.25
$ws.Range("K74").Value = 4356.25
$ws.Range("M74").Value = -3420.25
$ws.Range("H77").Value = 4759.533
$ws.Range("I77").Value = 4356.25
$ws.Range("K77").Value = 21781.25
$ws.Range("M77").Value = -17101.25
$ws.Range("H101").Value = 978.3333
$ws.Range("I101").Value = 450
$ws.Range("J101").Value = 1242.5
$ws.Range("K101").Value = 1350
$ws.Range("L101").Value = 3727.5
$ws.Range("M101").Value = 272
$ws.Range("N101").Value = -6971.5
$ws.Range("H102").Value = 34950
$ws.Range("J102").Value = 34950
$ws.Range("L102").Value = 34950
$ws.Range("N102").Value = -41440
$ws.Range("H132").Value = 17698.783
$ws.Range("I132").Value = 7475.522
$ws.Range("K132").Value = 22426.566
$ws.Range("M132").Value = -19896.566

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 2800
$ws.Range("I5").Value = 2800
$ws.Range("K5").Value = 2800
$ws.Range("M5").Value = -2688
$ws.Range("H32").Value = 13357.592
$ws.Range("I32").Value = 13497.136
$ws.Range("K32").Value = 13497.136
$ws.Range("M32").Value = -13210.136
$ws.Range("H97").Value = 478.69232
$ws.Range("I97").Value = 486.44446
$ws.Range("J97").Value = 461.25
$ws.Range("K97").Value = 486.44446
$ws.Range("L97").Value = 461.25
$ws.Range("M97").Value = 9.555540000000008
$ws.Range("N97").Value = -1453.25
$ws.Range("H102").Value = 361938.12
$ws.Range("I102").Value = 457849.34
$ws.Range("K102").Value = 457849.34
$ws.Range("M102").Value = -456227.34
$ws.Range("H132").Value = 25227.84
$ws.Range("I132").Value = 33466.5
$ws.Range("K132").Value = 100399.5
$ws.Range("M132").Value = -97869.5

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 2800
$ws.Range("I4").Value = 2800
$ws.Range("K4").Value = 2800
$ws.Range("M4").Value = -2685
$ws.Range("H28").Value = 40000
$ws.Range("J28").Value = 40000
$ws.Range("L28").Value = 40000
$ws.Range("N28").Value = -40588
$ws.Range("H58").Value = 28890
$ws.Range("I58").Value = 38000
$ws.Range("J58").Value = 19780
$ws.Range("K58").Value = 38000
$ws.Range("L58").Value = 19780
$ws.Range("M58").Value = -37706
$ws.Range("N58").Value = -20368
$ws.Range("H64").Value = 4630194.5
$ws.Range("I64").Value = 10417114
$ws.Range("J64").Value = 658.6
$ws.Range("K64").Value = 10417114
$ws.Range("L64").Value = 658.6
$ws.Range("M64").Value = -10416889
$ws.Range("N64").Value = -1108.6
$ws.Range("H67").Value = 4630194.5
$ws.Range("I67").Value = 10417114
$ws.Range("J67").Value = 658.6
$ws.Range("K67").Value = 10417114
$ws.Range("L67").Value = 658.6
$ws.Range("M67").Value = -10416334
$ws.Range("N67").Value = -2218.6
$ws.Range("H80").Value = 803
$ws.Range("I80").Value = 1270
$ws.Range("J80").Value = 709.6
$ws.Range("K80").Value = 1270
$ws.Range("L80").Value = 709.6
$ws.Range("M80").Value = -272
$ws.Range("N80").Value = -2705.6
$ws.Range("H83").Value = 803
$ws.Range("I83").Value = 1270
$ws.Range("J83").Value = 709.6
$ws.Range("K83").Value = 6350
$ws.Range("L83").Value = 3548
$ws.Range("M83").Value = -1358
$ws.Range("N83").Value = -13532
$ws.Range("H86").Value = 1861.25
$ws.Range("I86").Value = 2222.5
$ws.Range("J86").Value = 1500
$ws.Range("K86").Value = 2222.5
$ws.Range("L86").Value = 1500
$ws.Range("M86").Value = -1099.5
$ws.Range("N86").Value = -3746
$ws.Range("H89").Value = 1861.25
$ws.Range("I89").Value = 2222.5
$ws.Range("J89").Value = 1500
$ws.Range("K89").Value = 11112.5
$ws.Range("L89").Value = 7500
$ws.Range("M89").Value = -5496.5
$ws.Range("N89").Value = -18732
$ws.Range("H94").Value = 527750.4
$ws.Range("I94").Value = 1142588.6
$ws.Range("J94").Value = 746.1429000000001
$ws.Range("K94").Value = 1142588.6
$ws.Range("L94").Value = 746.1429000000001
$ws.Range("M94").Value = -1142137.6
$ws.Range("N94").Value = -1648.1429

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2265
$ws.Range("I16").Value = 2231.375
$ws.Range("K16").Value = 2231.375
$ws.Range("M16").Value = -1944.375
$ws.Range("H22").Value = 442.94736
$ws.Range("I22").Value = 246.63637
$ws.Range("J22").Value = 712.875
$ws.Range("K22").Value = 246.63637
$ws.Range("L22").Value = 712.875
$ws.Range("M22").Value = 103.36363
$ws.Range("N22").Value = -1412.875
$ws.Range("H113").Value = 2265
$ws.Range("I113").Value = 2231.375
$ws.Range("K113").Value = 2231.375
$ws.Range("M113").Value = -61.375
$ws.Range("H132").Value = 13899436
$ws.Range("I132").Value = 15884748
$ws.Range("K132").Value = 47654244
$ws.Range("M132").Value = -47651714
$ws.Range("H134").Value = 2758.1667
$ws.Range("I134").Value = 2758.1667
$ws.Range("K134").Value = 8274.500100000001
$ws.Range("M134").Value = -5739.500100000001

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 23383
$ws.Range("I120").Value = 5957.5
$ws.Range("J120").Value = 35000
$ws.Range("K120").Value = 17872.5
$ws.Range("L120").Value = 105000
$ws.Range("M120").Value = -13034.5
$ws.Range("N120").Value = -114676

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").ClearContents()
$ws.Range("N29").Value = 0
$ws.Range("H80").Value = 1122265.6
$ws.Range("I80").Value = 2772749.8
$ws.Range("J80").Value = 21942.889
$ws.Range("K80").Value = 2772749.8
$ws.Range("L80").Value = 21942.889
$ws.Range("M80").Value = -2771751.8
$ws.Range("N80").Value = -23938.889
$ws.Range("H83").Value = 1122265.6
$ws.Range("I83").Value = 2772749.8
$ws.Range("J83").Value = 21942.889
$ws.Range("K83").Value = 13863749
$ws.Range("L83").Value = 109714.445
$ws.Range("M83").Value = -13858757
$ws.Range("N83").Value = -119698.445
$ws.Range("H113").Value = 1526.6666
$ws.Range("J113").Value = 1706.25
$ws.Range("L113").Value = 1706.25
$ws.Range("N113").Value = -6046.25
$ws.Range("H132").Value = 574956.8
$ws.Range("I132").Value = 169992.25
$ws.Range("J132").Value = 1114909.5
$ws.Range("K132").Value = 509976.75
$ws.Range("L132").Value = 3344728.5
$ws.Range("M132").Value = -507446.75
$ws.Range("N132").Value = -3349788.5

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 7144747
$ws.Range("I16").Value = 11112624
$ws.Range("J16").Value = 2569.2
$ws.Range("K16").Value = 11112624
$ws.Range("L16").Value = 2569.2
$ws.Range("M16").Value = -11112454
$ws.Range("N16").Value = -2909.2
$ws.Range("H22").Value = 1250
$ws.Range("H27").Value = 1250
$ws.Range("H82").Value = 15625750
$ws.Range("I82").Value = 31250000
$ws.Range("K82").Value = 31250000
$ws.Range("M82").Value = -31249639
$ws.Range("H85").Value = 15625750
$ws.Range("I85").Value = 31250000
$ws.Range("K85").Value = 31250000
$ws.Range("M85").Value = -31248752

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 24000
$ws.Range("I21").Value = 24000
$ws.Range("K21").Value = 24000
$ws.Range("M21").Value = -23765
$ws.Range("H24").Value = 3000
$ws.Range("I24").Value = 3000
$ws.Range("K24").Value = 3000
$ws.Range("M24").Value = -2770
$ws.Range("H35").Value = 24000
$ws.Range("I35").Value = 24000
$ws.Range("K35").Value = 24000
$ws.Range("M35").Value = -23710
$ws.Range("H48").Value = 22500
$ws.Range("J48").Value = 22500
$ws.Range("L48").Value = 22500
$ws.Range("N48").Value = -23638
$ws.Range("H81").Value = 3212506.2
$ws.Range("I81").Value = 3479880.2
$ws.Range("J81").Value = 2983328.5
$ws.Range("K81").Value = 6959760.4
$ws.Range("L81").Value = 5966657
$ws.Range("M81").Value = -6958699.4
$ws.Range("N81").Value = -5968779
$ws.Range("H84").Value = 3212506.2
$ws.Range("I84").Value = 3479880.2
$ws.Range("J84").Value = 2983328.5
$ws.Range("K84").Value = 34798802
$ws.Range("L84").Value = 29833285
$ws.Range("M84").Value = -34793498
$ws.Range("N84").Value = -29843893
$ws.Range("H132").Value = 39683124
$ws.Range("I132").Value = 4274132
$ws.Range("J132").Value = 500000000
$ws.Range("K132").Value = 12822396
$ws.Range("L132").Value = 1500000000
$ws.Range("M132").Value = -12819866
$ws.Range("N132").Value = -1500005060
